$d = $word.ActiveDocument

# Remove the existing hidden "_GoBack" bookmark from paragraph 1; a fresh
# one will be recreated (via literal XML) in the new second paragraph.
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks.Item("_GoBack").Delete()
}

# Insert a brand-new paragraph at the end of the document containing
# "the second time." followed by the _GoBack bookmark and a trailing
# space run.
$endRange = $d.Content
$endRange.Collapse(0)
$xml = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' + `
    '<w:r><w:t>the second time.</w:t></w:r>' + `
    '<w:bookmarkStart w:id="0" w:name="_GoBack"/>' + `
    '<w:bookmarkEnd w:id="0"/>' + `
    '<w:r><w:t xml:space="preserve"> </w:t></w:r>' + `
    '</w:p>'
[void]$endRange.InsertXML($xml)
